# Auto update Excel log
# Appends newly-captured sensor-log rows to the PIR, Humidity, Temperature
# and mmWave sheets (rows were generated by the monitoring system for 2026-01-28).
$wb = $excel.ActiveWorkbook

# --- PIR sheet: append new sensor log rows ---
$ws = $wb.Worksheets.Item('PIR')
$newRows = @(
    @(34,'2026-01-28','16:07:52','16:00','Bathroom','No Motion','Inactive'),
    @(35,'2026-01-28','16:07:53','16:00','Bathroom','No Motion','Inactive'),
    @(36,'2026-01-28','16:07:56','16:00','Bathroom','No Motion','Inactive'),
    @(37,'2026-01-28','16:08:01','16:00','Bathroom','No Motion','Inactive'),
    @(38,'2026-01-28','16:08:06','16:00','Bathroom','No Motion','Inactive'),
    @(39,'2026-01-28','16:08:11','16:00','Bathroom','No Motion','Inactive'),
    @(40,'2026-01-28','16:08:16','16:00','Bathroom','No Motion','Inactive'),
    @(41,'2026-01-28','16:08:21','16:00','Bathroom','No Motion','Inactive'),
    @(42,'2026-01-28','16:08:26','16:00','Bathroom','No Motion','Inactive'),
    @(43,'2026-01-28','16:08:32','16:00','Bathroom','No Motion','Inactive'),
    @(44,'2026-01-28','16:08:37','16:00','Bathroom','No Motion','Inactive'),
    @(45,'2026-01-28','16:08:42','16:00','Bathroom','No Motion','Inactive'),
    @(46,'2026-01-28','16:08:47','16:00','Bathroom','No Motion','Inactive'),
    @(47,'2026-01-28','16:08:52','16:00','Bathroom','No Motion','Inactive')
)
foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# --- Humidity sheet: append new sensor log rows ---
$ws = $wb.Worksheets.Item('Humidity')
$newRows = @(
    @(34,'2026-01-28','16:07:52','16:00','Bathroom','88.7%','Active'),
    @(35,'2026-01-28','16:07:53','16:00','Bathroom','87.7%','Active'),
    @(36,'2026-01-28','16:07:56','16:00','Bathroom','88.7%','Active'),
    @(37,'2026-01-28','16:08:04','16:00','Bathroom','88.7%','Active'),
    @(38,'2026-01-28','16:08:08','16:00','Bathroom','88.6%','Active'),
    @(39,'2026-01-28','16:08:12','16:00','Bathroom','87.7%','Active'),
    @(40,'2026-01-28','16:08:16','16:00','Bathroom','88.7%','Active'),
    @(41,'2026-01-28','16:08:24','16:00','Bathroom','87.7%','Active'),
    @(42,'2026-01-28','16:08:28','16:00','Bathroom','88.8%','Active'),
    @(43,'2026-01-28','16:08:32','16:00','Bathroom','87.7%','Active'),
    @(44,'2026-01-28','16:08:36','16:00','Bathroom','87.2%','Active'),
    @(45,'2026-01-28','16:08:44','16:00','Bathroom','87.8%','Active'),
    @(46,'2026-01-28','16:08:48','16:00','Bathroom','88.7%','Active'),
    @(47,'2026-01-28','16:08:52','16:00','Bathroom','87.7%','Active')
)
foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# --- Temperature sheet: append new sensor log rows ---
$ws = $wb.Worksheets.Item('Temperature')
$newRows = @(
    @(34,'2026-01-28','16:07:53','16:00','Bathroom','22.7C','Active'),
    @(35,'2026-01-28','16:07:53','16:00','Bathroom','22.7C','Active'),
    @(36,'2026-01-28','16:07:56','16:00','Bathroom','22.7C','Active'),
    @(37,'2026-01-28','16:08:04','16:00','Bathroom','22.7C','Active'),
    @(38,'2026-01-28','16:08:08','16:00','Bathroom','22.7C','Active'),
    @(39,'2026-01-28','16:08:12','16:00','Bathroom','22.7C','Active'),
    @(40,'2026-01-28','16:08:16','16:00','Bathroom','22.7C','Active'),
    @(41,'2026-01-28','16:08:24','16:00','Bathroom','22.7C','Active'),
    @(42,'2026-01-28','16:08:28','16:00','Bathroom','22.8C','Active'),
    @(43,'2026-01-28','16:08:32','16:00','Bathroom','22.7C','Active'),
    @(44,'2026-01-28','16:08:36','16:00','Bathroom','22.7C','Active'),
    @(45,'2026-01-28','16:08:44','16:00','Bathroom','22.7C','Active'),
    @(46,'2026-01-28','16:08:48','16:00','Bathroom','22.8C','Active'),
    @(47,'2026-01-28','16:08:52','16:00','Bathroom','22.7C','Active')
)
foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# --- mmWave sheet: append new sensor log rows ---
$ws = $wb.Worksheets.Item('mmWave')
$newRows = @(
    @(9,'2026-01-28','16:08:14','16:00','Living Room','No Presence','Inactive'),
    @(10,'2026-01-28','16:08:19','16:00','Living Room','Presence Detected','Active'),
    @(11,'2026-01-28','16:08:25','16:00','Living Room','No Presence','Inactive'),
    @(12,'2026-01-28','16:08:49','16:00','Living Room','Presence Detected','Active')
)
foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}
